$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (rows 2, 3, 4) ---
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 17.25

# --- New year columns N (14) and O (15) on row 4 (headers 2022 / 2023) ---
# Copy formatting from the existing M4 (2021) cell, then set the new values.
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)
$ws.Cells.Item(4, 14).Value = 2022
$ws.Cells.Item(4, 15).Value = 2023

# --- New data columns N (14) and O (15) on row 5 (values 6.53 / 6.53) ---
# Copy formatting from the existing M5 (6.53) cell, then set the new values.
$ws.Range("M5").Copy()
$ws.Range("N5:O5").PasteSpecial(-4122)
$ws.Cells.Item(5, 14).Value = 6.53
$ws.Cells.Item(5, 15).Value = 6.53

$excel.CutCopyMode = 0
